$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''278.48'
$ws.Range("E2").Value = '''6.58%'
$ws.Range("G2").Value = '''5'
$ws.Range("D3").Value = '''27.26'
$ws.Range("E3").Value = '''0.78%'
$ws.Range("G3").Value = '''5'
$ws.Range("D4").Value = '''4.821'
$ws.Range("E4").Value = '''2.57%'
$ws.Range("G4").Value = '''5'
$ws.Range("D5").Value = '''0.06260'
$ws.Range("E5").Value = '''0.62%'
$ws.Range("G5").Value = '''5'
$ws.Range("D6").Value = '''6.864'
$ws.Range("E6").Value = '''1.77%'
$ws.Range("G6").Value = '''5'
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = '''0.8784'
$ws.Range("E7").Value = '''2.93%'
$ws.Range("G7").Value = '''5'
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = '''0.9435'
$ws.Range("E8").Value = '''3.29%'
$ws.Range("G8").Value = '''5'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '''0.1450'
$ws.Range("E9").Value = '''3.37%'
$ws.Range("G9").Value = '''5'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.05157'
$ws.Range("E10").Value = '''6.65%'
$ws.Range("G10").Value = '''5'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.07277'
$ws.Range("E11").Value = '''2.70%'
$ws.Range("G11").Value = '''5'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03136'
$ws.Range("E12").Value = '''0.98%'
$ws.Range("G12").Value = '''5'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09049'
$ws.Range("E13").Value = '''-0.20%'
$ws.Range("G13").Value = '''5'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001560'
$ws.Range("E14").Value = '''2.06%'
$ws.Range("G14").Value = '''5'
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = '''0.0006270'
$ws.Range("E15").Value = '''1.90%'
$ws.Range("G15").Value = '''5'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.005934'
$ws.Range("E16").Value = '''-2.69%'
$ws.Range("G16").Value = '''5'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.450'
$ws.Range("E17").Value = '''0.23%'
$ws.Range("G17").Value = '''5'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '''3.276'
$ws.Range("E18").Value = '''3.16%'
$ws.Range("G18").Value = '''5'
$ws.Range("E19").Value = '''4.64%'
$ws.Range("G19").Value = '''5'
$ws.Range("E20").Value = '''-0.61%'
$ws.Range("G20").Value = '''5'
$ws.Range("D21").Value = '''0.1310'
$ws.Range("E21").Value = '''-0.06%'
$ws.Range("G21").Value = '''5'
$ws.Range("D22").Value = '''3.845'
$ws.Range("E22").Value = '''-5.96%'
$ws.Range("G22").Value = '''5'
$ws.Range("D23").Value = '''0.04317'
$ws.Range("E23").Value = '''1.73%'
$ws.Range("G23").Value = '''5'
$ws.Range("D24").Value = '''0.001172'
$ws.Range("E24").Value = '''-3.44%'
$ws.Range("G24").Value = '''5'
$ws.Range("D25").Value = '''0.004276'
$ws.Range("E25").Value = '''4.63%'
$ws.Range("G25").Value = '''5'
$ws.Range("E26").Value = '''-0.18%'
$ws.Range("G26").Value = '''5'
$ws.Range("E27").Value = '''8.17%'
$ws.Range("G27").Value = '''5'
$ws.Range("G28").Value = '''5'
$ws.Range("G29").Value = '''5'
$ws.Range("G30").Value = '''5'
$ws.Range("G31").Value = '''5'
$ws.Range("G32").Value = '''5'
$ws.Range("G33").Value = '''5'
$ws.Range("G34").Value = '''5'
$ws.Range("G35").Value = '''5'
$ws.Range("G36").Value = '''5'
$ws.Range("G37").Value = '''5'
$ws.Range("G38").Value = '''5'
$ws.Range("G39").Value = '''5'
$ws.Range("D40").Value = '''0.04031'
$ws.Range("E40").Value = '''2.92%'
$ws.Range("G40").Value = '''5'
$ws.Range("D41").Value = '''0.006482'
$ws.Range("E41").Value = '''57.46%'
$ws.Range("G41").Value = '''5'
$ws.Range("D42").Value = '''0.1153'
$ws.Range("E42").Value = '''3.64%'
$ws.Range("G42").Value = '''5'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '''0.002157'
$ws.Range("E43").Value = '''-2.43%'
$ws.Range("G43").Value = '''5'
$ws.Range("B44").Value = 'LocalTraders'
$ws.Range("C44").Value = 'https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct'
$ws.Range("D44").Value = '''0.01405'
$ws.Range("E44").Value = '''1.10%'
$ws.Range("G44").Value = '''5'
$ws.Range("D45").Value = '''0.00005153'
$ws.Range("E45").Value = '''-0.18%'
$ws.Range("G45").Value = '''5'
$ws.Range("E46").Value = '''-0.18%'
$ws.Range("G46").Value = '''5'
$ws.Range("D47").Value = '''2.346'
$ws.Range("E47").Value = '''665.70%'
$ws.Range("G47").Value = '''5'
$ws.Range("G48").Value = '''5'
$ws.Range("E49").Value = '''-0.18%'
$ws.Range("G49").Value = '''5'
$ws.Range("E50").Value = '''-0.18%'
$ws.Range("G50").Value = '''5'
$ws.Range("G51").Value = '''5'

Write-Host "Applied 143 cell updates"
